$wb = $excel.ActiveWorkbook

# --- Rename sheets ---
# "ACD_Data" (currently the active/selected tab, r:id rId1) becomes "ACD_Data_bob".
# "ACD_Data_old" (r:id rId2) becomes "ACD_Data" and becomes the new active/selected tab.
$wsNew = $wb.Worksheets.Item("ACD_Data")
$wsOld = $wb.Worksheets.Item("ACD_Data_old")

$wsNew.Name = "ACD_Data_bob"
$wsOld.Name = "ACD_Data"

# Make the renamed "ACD_Data" sheet the active / selected tab (moves tabSelected
# from the old ACD_Data_bob sheet to this one, and updates the workbook's
# active tab index accordingly).
$wsOld.Select()

# --- Adjust row height on Data_Dictionary sheet ---
# Row 20 (the "Tail_Height_at_OEW_ft" entry) grows from 31pt to 46.5pt.
$wsDict = $wb.Worksheets.Item("Data_Dictionary")
$wsDict.Rows.Item(20).RowHeight = 46.5
